$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 59, shifting rows 59:111 down to 60:112
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with data (copy of old row 59 values,
# but with updated date / volume / price fields reflecting the new weekly record)
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 45096
$ws.Range("D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = 300000001
$ws.Range("G59").Value = "Rabanito"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 65
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 8000
$ws.Range("M59").Value = 8000
$ws.Range("N59").Value = "$/docena de paquetes"
$ws.Range("O59").Value = "Provincia de Cautín"
$ws.Range("P59").Value = 667
$ws.Range("Q59").Value = 12
$ws.Range("R59").Value = "Hortaliza"
